$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column A to fit the new "Khushboo / Aman" label
$ws.Columns.Item(1).ColumnWidth = 22.44140625

# Row 14 already exists: update the name + notes text to reflect the added collaborator
$ws.Range("A14").Value = "Khushboo / Aman"
$ws.Range("H14").Value = "Making login logout and room list responsive with Aman"

# New row 15: dark mode theme work
$ws.Range("A15").Value = "Khushboo / Aman"
$ws.Range("C15").Value = (Get-Date -Year 2020 -Month 8 -Day 7)
$ws.Range("E15").Value = 250
$ws.Range("G15").Value = "https://github.com/chat-loc/chatloc.github.io/commit/286ca9c3a8555e0c7234327dfac620043f4d2921"
$ws.Hyperlinks.Add($ws.Range("G15"), "https://github.com/chat-loc/chatloc.github.io/commit/286ca9c3a8555e0c7234327dfac620043f4d2921", "", "", "https://github.com/chat-loc/chatloc.github.io/commit/286ca9c3a8555e0c7234327dfac620043f4d2921") | Out-Null
$ws.Range("H15").Value = "dark mode theme with Aman"

# New row 16: merge to master
$ws.Range("A16").Value = "Khushboo / Aman"
$ws.Range("C16").Value = (Get-Date -Year 2020 -Month 8 -Day 8)
$ws.Range("E16").Value = 10
$ws.Range("G16").Value = "https://github.com/chat-loc/chatloc.github.io/commit/fa82cd20851166d39e7316799dab1363fe4a4404"
$ws.Hyperlinks.Add($ws.Range("G16"), "https://github.com/chat-loc/chatloc.github.io/commit/fa82cd20851166d39e7316799dab1363fe4a4404", "", "", "https://github.com/chat-loc/chatloc.github.io/commit/fa82cd20851166d39e7316799dab1363fe4a4404") | Out-Null
$ws.Range("H16").Value = "merging the dark mode theme to master branch with Aman"

# Match styles of the neighboring rows for the new data rows
$ws.Range("A15:A16").Style = $ws.Range("A14").Style
$ws.Range("C15:C16").Style = $ws.Range("C14").Style
$ws.Range("E15:E16").Style = $ws.Range("E14").Style
$ws.Range("G15:G16").Style = $ws.Range("G14").Style
$ws.Range("H15:H16").Style = $ws.Range("H14").Style

# Update the active selection like the author's last-saved state
$ws.Range("E17").Select()
